$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested count) values on both the "展览" and
# "全部类型" worksheets, which carry duplicate copies of the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F12").Value = 965
    $ws.Range("F14").Value = 539
}
